{"js": "// Add the three GaN character styles (GaNStyle, GaNParagraph, GaNLinks) and\n// apply GaNParagraph / GaNLinks to the relevant runs, matching the commit\n// \"Add styles to the new paragraphs\".\n\n// --- helper: create a character style and set its font props, re-acquiring\n// the style through the styles collection so the font writes stick. ---\nasync function addCharacterStyle(context, name, fontProps) {\n  context.document.addStyle(name, Word.StyleType.character);\n  await context.sync();\n\n  const styles = context.document.getStyles();\n  styles.load(\"items/nameLocal\");\n  await context.sync();\n\n  let target = null;\n  for (let i = 0; i < styles.items.length; i++) {\n    if (styles.items[i].nameLocal === name) {\n      target = styles.items[i];\n      break;\n    }\n  }\n  if (target) {\n    if (fontProps.name !== undefined) target.font.name = fontProps.name;\n    if (fontProps.size !== undefined) target.font.size = fontProps.size;\n    if (fontProps.bold !== undefined) target.font.bold = fontProps.bold;\n    if (fontProps.color !== undefined) target.font.color = fontProps.color;\n    if (fontProps.underline !== undefined) target.font.underline = fontProps.underline;\n    await context.sync();\n  }\n  return target;\n}\n\n// 1) GaNStyle \u2014 Calibri, 14pt (sz 28), no other attributes.\nawait addCharacterStyle(context, \"GaNStyle\", { name: \"Calibri\", size: 14 });\n\n// 2) GaNParagraph \u2014 Calibri, 10pt (sz 20).\nawait addCharacterStyle(context, \"GaNParagraph\", { name: \"Calibri\", size: 10 });\n\n// 3) GaNLinks \u2014 Calibri, bold, navy (000080), 9.5pt (sz 19), single underline.\nawait addCharacterStyle(context, \"GaNLinks\", {\n  name: \"Calibri\",\n  bold: true,\n  color: \"#000080\",\n  size: 9.5,\n  underline: Word.UnderlineType.single,\n});\n\n// --- apply GaNParagraph to every run with the recurring Swedish paragraph ---\nconst paragraphText =\n  \"Du deltar i en v\u00e4rldsomsp\u00e4nnande kampanj f\u00f6r att observera och rapportera de svagaste synliga stj\u00e4rnorna, som ett m\u00e5tt p\u00e5 ljusf\u00f6roreningarna p\u00e5 orten. Genom att hitta och observera Cygnus konstellation p\u00e5 natthimlen kan folk i hela v\u00e4rlden l\u00e4ra sig hur belysningen i v\u00e5ra samh\u00e4llen och omgivningar bidrar till ljusf\u00f6roreningar. Era bidrag till online-databasen hj\u00e4lper till att dokumentera den synliga natthimlens \u00f6ver hela v\u00e4rlden.\";\n\nconst paragraphHits = context.document.body.search(paragraphText, { matchCase: true });\nparagraphHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphHits.items.length; i++) {\n  paragraphHits.items[i].style = \"GaNParagraph\";\n}\nawait context.sync();\n\n// --- apply GaNLinks to the run holding the GaNight link text ---\nconst linkText = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\n\nconst linkHits = context.document.body.search(linkText, { matchCase: true });\nlinkHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < linkHits.items.length; i++) {\n  linkHits.items[i].style = \"GaNLinks\";\n}\nawait context.sync();\n", "ps1": "# Add the three GaN character styles (GaNStyle, GaNParagraph, GaNLinks) and\n# apply GaNParagraph / GaNLinks to the relevant runs, matching the commit\n# \"Add styles to the new paragraphs\".\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeCharacter = 2 (named constants aren't pre-seeded in this host)\n$wdStyleTypeCharacter = 2\n$wdUnderlineSingle = 1\n\n# 1) GaNStyle - Calibri, 14pt (sz 28), no other attributes.\n$gaNStyle = $d.Styles.Add(\"GaNStyle\", $wdStyleTypeCharacter)\n$gaNStyle.Font.Name = \"Calibri\"\n$gaNStyle.Font.Size = 14\n\n# 2) GaNParagraph - Calibri, 10pt (sz 20).\n$gaNParagraph = $d.Styles.Add(\"GaNParagraph\", $wdStyleTypeCharacter)\n$gaNParagraph.Font.Name = \"Calibri\"\n$gaNParagraph.Font.Size = 10\n\n# 3) GaNLinks - Calibri, bold, navy (000080 = RGB(0,0,128) = 8388608), 9.5pt (sz 19), single underline.\n$gaNLinks = $d.Styles.Add(\"GaNLinks\", $wdStyleTypeCharacter)\n$gaNLinks.Font.Name = \"Calibri\"\n$gaNLinks.Font.Bold = $true\n$gaNLinks.Font.Color = 8388608\n$gaNLinks.Font.Size = 9.5\n$gaNLinks.Font.Underline = $wdUnderlineSingle\n\n# --- apply GaNParagraph to every run with the recurring Swedish paragraph ---\n$paragraphText = \"Du deltar i en v\u00e4rldsomsp\u00e4nnande kampanj f\u00f6r att observera och rapportera de svagaste synliga stj\u00e4rnorna, som ett m\u00e5tt p\u00e5 ljusf\u00f6roreningarna p\u00e5 orten. Genom att hitta och observera Cygnus konstellation p\u00e5 natthimlen kan folk i hela v\u00e4rlden l\u00e4ra sig hur belysningen i v\u00e5ra samh\u00e4llen och omgivningar bidrar till ljusf\u00f6roreningar. Era bidrag till online-databasen hj\u00e4lper till att dokumentera den synliga natthimlens \u00f6ver hela v\u00e4rlden.\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $paragraphText\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\nwhile ($rng.Find.Execute()) {\n  $rng.Style = \"GaNParagraph\"\n  $rng.Collapse(0)\n}\n\n# --- apply GaNLinks to the run holding the GaNight link text ---\n$linkText = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = $linkText\n$rng2.Find.Forward = $true\n$rng2.Find.Wrap = 0\nwhile ($rng2.Find.Execute()) {\n  $rng2.Style = \"GaNLinks\"\n  $rng2.Collapse(0)\n}\n"}
